$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Add the new data row (26 April 2020) right after the last existing row.
# Using the worksheet Table (ListObject) so the table range auto-expands.
$table = $ws.ListObjects.Item("Table3")
$newRow = $table.ListRows.Add()

$ws.Range("A46").Value = 43947
$ws.Range("B46").Value = 30177
$ws.Range("C46").Value = 2357
$ws.Range("D46").Value = 99
$ws.Range("E46").Value = 3558

# Update the view: scroll so row 33 is the top row and select B47 (just below the new data).
$ws.Range("B47").Select()
$excel.ActiveWindow.ScrollRow = 33
